$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176199555397034
$ws.Range("B1").Value = 2.413364171981812
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.337907314300537
$ws.Range("E1").Value = 1.202196359634399
